# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on row 2 of the
# per-language report sheets (zh-cn, de-de) to reflect a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 02:37:32"
$wsZhCn.Range("H2").Value = "2016-03-19 02:37:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 02:37:36"
$wsDeDe.Range("H2").Value = "2016-03-19 02:37:57"
